$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "ALLDO followed by 4 values..." -> "DSSP followed by 4 values..."
#    Word's "_GoBack" last-edit bookmark moves to sit right after the
#    newly typed replacement text, splitting that run in two.
# -----------------------------------------------------------------

# The _GoBack bookmark currently sits after "HOME defaults to 90."
# (paragraph 22). Remove it from there first.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# Find the "ALLDO" paragraph and replace the word with "DSSP".
$alldoPara = $d.Paragraphs(18)
$alldoStart = $alldoPara.Range.Start
$alldoWord = $d.Range($alldoStart, $alldoStart + 5)
$alldoWord.Text = "DSSP"

# Re-create _GoBack immediately after the replacement text, which
# splits the run and matches Word's own "last edit" bookmark tracking.
$newGoBackPos = $d.Range($alldoStart + 4, $alldoStart + 4)
$d.Bookmarks.Add("_GoBack", $newGoBackPos)

# -----------------------------------------------------------------
# 2) Add a new bold closing line: Device is called "EarGear" over
#    Bluetooth, in the last (previously empty) paragraph.
# -----------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range

# Apply bold formatting to the still-empty paragraph first so the
# paragraph mark's own run properties (pPr/rPr) pick up both the
# western and complex-script bold flags.
$lastRange.Font.Name = "Century Gothic"
$lastRange.Bold = 1
$lastRange.BoldBi = 1
$lastRange.Text = [char]0x201C + "EarGear" + [char]0x201D
$lastRange.InsertBefore("Device is called ")
$lastRange.InsertAfter(" over Bluetooth")

# Re-apply bold formatting across the now-populated paragraph so the
# inserted runs also carry the bold / complex-script-bold flags.
$finalRange = $lastPara.Range
$finalRange.Font.Name = "Century Gothic"
$finalRange.Bold = 1
$finalRange.BoldBi = 1
